# Add the "2022-Q1" data sheet and roll its totals into "总计".

$wb = $excel.ActiveWorkbook

# Use the "2021-Q4" sheet as a style/layout template: it already has the
# exact header row (基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名)
# with the bold+bordered "s=2" style, and column-A uses the same style too.
$template = $wb.Worksheets.Item("2021-Q4")

# Create the new sheet. It gets inserted at the front, so (re)fetch "总计"'s
# reference only *after* that insertion shifted everyone's index, then move
# the new sheet to sit right before it (tab order becomes 2020-Q4, 2021-Q3,
# 2021-Q4, 2022-Q1, 总计).
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"
$total = $wb.Worksheets.Item("总计")
$newSheet.Move($total)

# Fetch a fresh reference by name: after Move() the old variable keeps
# pointing at the sheet that now occupies the old slot, not the moved one.
$ws = $wb.Worksheets.Item("2022-Q1")
$total = $wb.Worksheets.Item("总计")

# Copy header formatting/text and the column-A style from the template sheet.
$template.Range("B1:H1").Copy($ws.Range("B1:H1"))
$template.Range("A2").Copy($ws.Range("A2"))
$template.Range("A2").Copy($ws.Range("A3"))
$template.Range("A2").Copy($ws.Range("A4"))
$template.Range("A2").Copy($ws.Range("A5"))
$template.Range("A2").Copy($ws.Range("A6"))

# The fund-code / numeric-looking columns are stored as text, not numbers,
# so force a text number format before assigning them.
$ws.Range("B2:B6").NumberFormat = "@"
$ws.Range("D2:G6").NumberFormat = "@"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "501007"
$ws.Range("C2").Value = "汇添富中证互联网医疗主题指数（LOF）A"
$ws.Range("D2").Value = "0.58"
$ws.Range("E2").Value = "93.89"
$ws.Range("F2").Value = "5.23"
$ws.Range("G2").Value = "0.0303"
$ws.Range("H2").Value = 3

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "168701"
$ws.Range("C3").Value = "合煦智远国证香蜜湖金融科技指数(LOF)A"
$ws.Range("D3").Value = "0.90"
$ws.Range("E3").Value = "93.15"
$ws.Range("F3").Value = "1.66"
$ws.Range("G3").Value = "0.0149"
$ws.Range("H3").Value = 10

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "501008"
$ws.Range("C4").Value = "汇添富中证互联网医疗主题指数（LOF）C"
$ws.Range("D4").Value = "0.19"
$ws.Range("E4").Value = "93.89"
$ws.Range("F4").Value = "5.23"
$ws.Range("G4").Value = "0.0099"
$ws.Range("H4").Value = 3

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "165522"
$ws.Range("C5").Value = "信诚中证TMT产业主题指数（LOF）"
$ws.Range("D5").Value = "0.58"
$ws.Range("E5").Value = "93.74"
$ws.Range("F5").Value = "1.14"
$ws.Range("G5").Value = "0.0066"
$ws.Range("H5").Value = 10

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "168702"
$ws.Range("C6").Value = "合煦智远国证香蜜湖金融科技指数(LOF)C"
$ws.Range("D6").Value = "0.22"
$ws.Range("E6").Value = "93.15"
$ws.Range("F6").Value = "1.66"
$ws.Range("G6").Value = "0.0037"
$ws.Range("H6").Value = 10

# Now roll the new quarter into the "总计" (totals) summary sheet: insert a
# fresh row right under the header and push the existing quarters down.
$total.Rows.Item(2).Insert()

# The inserted row picked up the header's bold/bordered formatting; clear
# the data cells back to plain formatting and restyle column A to match
# the other index cells (style "2": bold + thin border, from the template).
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 5
$total.Range("D2").Value = 0.07000000000000001

# Renumber the (0-based) index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
